$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (J2, K2, M2 are standalone formulas) ---
$ws.Range("J2").Formula = "=F2*H2 / (1-H2)"
$ws.Range("K2").Formula = "=ROUNDUP(J2*1000, 0)*5"
# M2 formula (=K2+G2) is unchanged; its cached value will be refreshed automatically.

# --- J3:J6 become one shared-formula group (fuel cost incl. extra weight) ---
$ws.Range("J3:J6").Formula = "=F3*H3 / (1-H3)"

# --- K3:K5 become a (narrower) shared-formula group using ROUNDUP ---
$ws.Range("K3:K5").Formula = "=ROUNDUP(J3*1000, 0)*5"

# K6 keeps its own original-style formula (J6*1000*5), recomputed with the new J6
$ws.Range("K6").Formula = "=J6*1000*5"

# M3:M5 formula text is unchanged (=K3+G3 shifted); values refresh automatically.
# M6 formula text is unchanged (=SUM(M2:M5)); value refreshes automatically.

# --- Row 8-11: min/max cost block ---
# J8, J10 and J11 form a shared-formula group; J9 keeps an explicit (unshared) formula
$ws.Range("J8:J11").Formula = "=(D2 + F2)*H2 / (1-H2)"
$ws.Range("J9").Formula = "=(D3 + F3)*H3 / (1-H3)"

# K8 becomes a standalone ROUNDUP-based formula
$ws.Range("K8").Formula = "=ROUNDUP(J8*1000, 0)*5"

# K9:K11 keep being a shared-formula group, now using ROUNDUP
$ws.Range("K9:K11").Formula = "=ROUNDUP(J9*1000, 0)*5"

# M8 now adds G2 directly instead of M2
$ws.Range("M8").Formula = "=K8+G2"

# M9:M11 become a shared-formula group (=K+G, shifted)
$ws.Range("M9:M11").Formula = "=K9+G3"

# K12 / M12 formulas (=SUM(...)) are unchanged; their cached values refresh automatically.

# --- Cosmetic: update the saved selection to M17, like the target workbook ---
$ws.Range("M17").Select()
